$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1 = 14, Q1 = 15 with same style as O1 (bold, centered/top, bordered)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Data rows 2-25: swap values in columns I,K,M,O and add new columns P,Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
